$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "PASS" status in column L for each data row (rows 2-8),
# matching the new STATUS column values recorded for each test case.
$ws.Range("L2:L8").Value = "PASS"
